$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 81: Volunteering with Staff | Hallowed Chestnut Wand
$ws.Range("H81").Value = 39950
$ws.Range("J81").Value = 39950
$ws.Range("L81").Value = 39950
$ws.Range("N81").Value = -41946

# Row 84: Scripture Is the Best Medicine (L) | Hallowed Chestnut Wand
$ws.Range("H84").Value = 39950
$ws.Range("J84").Value = 39950
$ws.Range("L84").Value = 119850
$ws.Range("N84").Value = -129834

# Row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 682.7895
$ws.Range("I98").Value = 709.3889
$ws.Range("J98").Value = 204
$ws.Range("K98").Value = 709.3889
$ws.Range("L98").Value = 204
$ws.Range("M98").Value = 788.6111
$ws.Range("N98").Value = -3200

# Row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 682.7895
$ws.Range("I122").Value = 709.3889
$ws.Range("J122").Value = 204
$ws.Range("K122").Value = 2128.1667
$ws.Range("L122").Value = 612
$ws.Range("M122").Value = 321.8332999999998
$ws.Range("N122").Value = -5512

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 315619.56
$ws.Range("I32").Value = 1237.08
$ws.Range("K32").Value = 1237.08
$ws.Range("M32").Value = -950.0799999999999

# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 1812
$ws.Range("I45").Value = 1540.2858
$ws.Range("J45").Value = 2192.4
$ws.Range("K45").Value = 1540.2858
$ws.Range("L45").Value = 2192.4
$ws.Range("M45").Value = -1163.2858
$ws.Range("N45").Value = -2946.4

# Row 98: Greaving | Doman Iron Greaves of Maiming
$ws.Range("H98").Value = 40500
$ws.Range("J98").Value = 40500
$ws.Range("L98").Value = 40500
$ws.Range("N98").Value = -46490

# Row 110: Scheduled Maintenance | Deepgold Ingot
$ws.Range("H110").Value = 1407.9412
$ws.Range("I110").Value = 913.05884
$ws.Range("J110").Value = 1902.8235
$ws.Range("K110").Value = 913.05884
$ws.Range("L110").Value = 1902.8235
$ws.Range("M110").Value = 1131.94116
$ws.Range("N110").Value = -5992.8235

# Row 133: Shielding My Students | Mountain Chromite Tower Shield
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 715316.5
$ws.Range("I86").Value = 762.7
$ws.Range("J86").Value = 2501701
$ws.Range("K86").Value = 762.7
$ws.Range("L86").Value = 2501701
$ws.Range("M86").Value = 360.3
$ws.Range("N86").Value = -2503947

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 715316.5
$ws.Range("I89").Value = 762.7
$ws.Range("J89").Value = 2501701
$ws.Range("K89").Value = 3813.5
$ws.Range("L89").Value = 12508505
$ws.Range("M89").Value = 1802.5
$ws.Range("N89").Value = -12519737

# Row 99: Meddle in Metal | Oroshigane Ingot
$ws.Range("H99").Value = 1189.0952
$ws.Range("I99").Value = 872.1579
$ws.Range("K99").Value = 872.1579
$ws.Range("M99").Value = 625.8421

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 4771.875
$ws.Range("I31").Value = 3290.8
$ws.Range("J31").Value = 7240.3335
$ws.Range("K31").Value = 3290.8
$ws.Range("L31").Value = 7240.3335
$ws.Range("M31").Value = -2995.8
$ws.Range("N31").Value = -7830.3335

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 4771.875
$ws.Range("I34").Value = 3290.8
$ws.Range("J34").Value = 7240.3335
$ws.Range("K34").Value = 3290.8
$ws.Range("L34").Value = 7240.3335
$ws.Range("M34").Value = -3088.8
$ws.Range("N34").Value = -7644.3335

# Row 54: The Turning Point | Garnet Grinding Wheel
$ws.Range("H54").Value = 23073.6
$ws.Range("J54").Value = 23073.6
$ws.Range("L54").Value = 23073.6
$ws.Range("N54").Value = -24389.6

# Row 75: The Darkest Hearth | Dark Chestnut Spinning Wheel
$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41996

# Row 78: Fruit of the Loom (L) | Dark Chestnut Spinning Wheel
$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -129984

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 1801.8889
$ws.Range("I134").Value = 1451.1428
$ws.Range("J134").Value = 3029.5
$ws.Range("K134").Value = 4353.428400000001
$ws.Range("L134").Value = 9088.5
$ws.Range("M134").Value = -1818.428400000001
$ws.Range("N134").Value = -14158.5

$ws = $wb.Worksheets.Item("CUL")
# Row 11: Putting the Squeeze On | Orange Juice
$ws.Range("H11").Value = 2078.7805
$ws.Range("I11").Value = 2439.182
$ws.Range("J11").Value = 1946.6333
$ws.Range("K11").Value = 7317.545999999999
$ws.Range("L11").Value = 5839.8999
$ws.Range("M11").Value = -7177.545999999999
$ws.Range("N11").Value = -6119.8999

# Row 57: The Egg Files | Deviled Eggs
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").Value = ""

# Row 96: Hunger Is No Game | Popoto Soba
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = ""

# Row 141: Ocean Explosion | Acqua Pazza
$ws.Range("H141").Value = 8020
$ws.Range("I141").Value = 8020
$ws.Range("K141").Value = 24060
$ws.Range("M141").Value = -18880

$ws = $wb.Worksheets.Item("GSM")
# Row 44: Actually, It's Loyalty | Aquamarine Bracelet
$ws.Range("H44").Value = 13214.286
$ws.Range("J44").Value = 9500
$ws.Range("L44").Value = 9500
$ws.Range("N44").Value = -10692

# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 2881.6667
$ws.Range("J102").Value = 5057.143
$ws.Range("L102").Value = 5057.143
$ws.Range("N102").Value = -8301.143

# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 3527.55
$ws.Range("I126").Value = 3071.7778
$ws.Range("K126").Value = 9215.3334
$ws.Range("M126").Value = -6745.3334

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore | Hard Leather
$ws.Range("H16").Value = 9913.571
$ws.Range("I16").Value = 9673.75
$ws.Range("J16").Value = 10233.333
$ws.Range("K16").Value = 9673.75
$ws.Range("L16").Value = 10233.333
$ws.Range("M16").Value = -9503.75
$ws.Range("N16").Value = -10573.333

# Row 22: Skin off Their Backs | Aldgoat Leather
$ws.Range("H22").Value = 911.5833
$ws.Range("J22").Value = 1048.75
$ws.Range("L22").Value = 1048.75
$ws.Range("N22").Value = -1638.75

# Row 27: Fire and Hide | Aldgoat Leather
$ws.Range("H27").Value = 911.5833
$ws.Range("J27").Value = 1048.75
$ws.Range("L27").Value = 1048.75
$ws.Range("N27").Value = -1262.75

# Row 68: You Could Say It's a Moving Target | Wyvern Leather
$ws.Range("H68").Value = 4384.615
$ws.Range("I68").Value = 4120
$ws.Range("J68").Value = 5266.6665
$ws.Range("K68").Value = 4120
$ws.Range("L68").Value = 5266.6665
$ws.Range("M68").Value = -3371
$ws.Range("N68").Value = -6764.6665

# Row 71: They Call It Bloody Mary (L) | Wyvern Leather
$ws.Range("H71").Value = 4384.615
$ws.Range("I71").Value = 4120
$ws.Range("J71").Value = 5266.6665
$ws.Range("K71").Value = 20600
$ws.Range("L71").Value = 26333.3325
$ws.Range("M71").Value = -16856
$ws.Range("N71").Value = -33821.3325

# Row 80: Don't Sweat the Small Fry | Dragonskin Wristbands
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""

# Row 83: It's All in the Wrists (L) | Dragonskin Wristbands
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""

# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 2782
$ws.Range("I122").Value = 2328.1667
$ws.Range("J122").Value = 5505
$ws.Range("K122").Value = 6984.500100000001
$ws.Range("L122").Value = 16515
$ws.Range("M122").Value = -4534.500100000001
$ws.Range("N122").Value = -21415

# Row 133: The Perfect Accessory | Loboskin Amulet of Fending
$ws.Range("H133").Value = 78994
$ws.Range("J133").Value = 78994
$ws.Range("L133").Value = 78994
$ws.Range("N133").Value = -84054

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke | Rainbow Cloth
$ws.Range("H62").Value = 5761.1763
$ws.Range("I62").Value = 2691.4285
$ws.Range("J62").Value = 7910
$ws.Range("K62").Value = 2691.4285
$ws.Range("L62").Value = 7910
$ws.Range("M62").Value = -2067.4285
$ws.Range("N62").Value = -9158

# Row 65: Desperate for Diversionaries (L) | Rainbow Cloth
$ws.Range("H65").Value = 5761.1763
$ws.Range("I65").Value = 2691.4285
$ws.Range("J65").Value = 7910
$ws.Range("K65").Value = 13457.1425
$ws.Range("L65").Value = 39550
$ws.Range("M65").Value = -10337.1425
$ws.Range("N65").Value = -45790

# Row 68: What Not to Wear | Holy Rainbow Shirt of Striking
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""

# Row 71: Appeal of Foreign Apparel (L) | Holy Rainbow Shirt of Striking
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 2310.1667
$ws.Range("I132").Value = 2302
$ws.Range("J132").Value = 2400
$ws.Range("K132").Value = 6906
$ws.Range("L132").Value = 7200
$ws.Range("M132").Value = -4376
$ws.Range("N132").Value = -12260
